# "fixed english memory game"
# - bump the balance on row 6 (F6: 1000 -> 1032)
# - append 37 new "played game" rows (77-113), each identical to the
#   existing template row 76 (moses / bro / 1234 / m@g.c / Male / 0)
#   so the new cells carry exactly the same types/styles as the rows
#   already produced by the game logger.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the balance for the row-6 user.
$ws.Range("F6").Value = 1032

# Duplicate the template row (row 76) into rows 77 through 113, just like
# the game appends a fresh "played" record per round, via copy/paste so the
# new cells pick up the identical shared-string typing/styling as the
# existing rows instead of being re-typed (which would turn the numeric
# looking "1234" string into a real number).
$template = $ws.Range("A76:F76")
for ($row = 77; $row -le 113; $row++) {
    $template.Copy()
    $dest = $ws.Range("A" + $row + ":F" + $row)
    $dest.PasteSpecial()
}
